$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.637.31"
$ws.Range("E2").Value = "  +1.43%  "

# Row 3
$ws.Range("D3").Value = "1.890.30"

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'238.88"
$ws.Range("E5").Value = "  +1.29%  "

# Row 6
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").Value = "'0.4830"

# Row 8
$ws.Range("D8").Value = "'0.2863"
$ws.Range("E8").Value = "  +2.15%  "

# Row 9
$ws.Range("D9").Value = "'0.06553"
$ws.Range("E9").Value = "  +1.27%  "

# Row 10
$ws.Range("D10").Value = "1.953.15"
$ws.Range("E10").Value = "  +5.34%  "

# Row 11
$ws.Range("D11").Value = "'0.07475"
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("D12").Value = "'16.69"
$ws.Range("E12").Value = "  +3.18%  "

# Row 13
$ws.Range("D13").Value = "'5.104"

# Row 14
$ws.Range("D14").Value = "'88.17"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15
$ws.Range("D15").Value = "'0.6679"
$ws.Range("E15").Value = "  +3.31%  "

# Row 16
$ws.Range("D16").Value = "30.610.67"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17
$ws.Range("D17").Value = "'13.28"
$ws.Range("E17").Value = "  +0.83%  "

# Row 18
$ws.Range("E18").Value = "  +0.10%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007577"
$ws.Range("E19").Value = "  -0.38%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'233.03"
$ws.Range("E20").Value = "  +4.09%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.127.33"
$ws.Range("E21").Value = "  +1.34%  "

# Row 22
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("D23").Value = "'5.283"
$ws.Range("E23").Value = "  +0.08%  "

# Row 24
$ws.Range("D24").Value = "'6.223"
$ws.Range("E24").Value = "  +2.58%  "

# Row 25
$ws.Range("D25").Value = "'169.36"
$ws.Range("E25").Value = "  +3.37%  "

# Row 26
$ws.Range("E26").Value = "  +1.56%  "

# Row 27
$ws.Range("D27").Value = "'18.83"
$ws.Range("E27").Value = "  +2.12%  "

# Row 28
$ws.Range("D28").Value = "'1.962"
$ws.Range("E28").Value = "  +2.18%  "

# Row 29
$ws.Range("D29").Value = "'0.1020"
$ws.Range("E29").Value = "  +11.10%  "

# Row 30
$ws.Range("D30").Value = "'1.398"
$ws.Range("E30").Value = "  -3.02%  "

# Row 31
$ws.Range("D31").Value = "'4.328"
$ws.Range("E31").Value = "  +1.96%  "

# Row 32
$ws.Range("D32").Value = "'4.030"
$ws.Range("E32").Value = "  +1.94%  "

# Row 33
$ws.Range("D33").Value = "'0.05068"
$ws.Range("E33").Value = "  +1.48%  "

# Row 34
$ws.Range("E34").Value = "  +5.99%  "

# Row 35
$ws.Range("D35").Value = "'0.7534"
$ws.Range("E35").Value = "  +2.53%  "

# Row 36
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.714"
$ws.Range("E37").Value = "  +0.86%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01875"
$ws.Range("E38").Value = "  +2.72%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.650"
$ws.Range("E39").Value = "  +1.74%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.9195"
$ws.Range("E40").Value = "  +2.27%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.073"
$ws.Range("E41").Value = "  +1.33%  "

# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'107.17"
$ws.Range("E42").Value = "  +0.85%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4295"
$ws.Range("E43").Value = "  +1.32%  "

# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'1.003"
$ws.Range("E44").Value = "  +0.23%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.646"
$ws.Range("E45").Value = "  -4.97%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.431"
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'64.29"
$ws.Range("E47").Value = "  +0.28%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1274"
$ws.Range("E48").Value = "  -3.34%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.497"
$ws.Range("E49").Value = "  -2.48%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.000"
$ws.Range("E50").Value = "  +2.77%  "

# Row 51
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'33.95"
$ws.Range("E51").Value = "  -0.11%  "
